$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.181.73"
$ws.Range("E2").Value = "  -0.79%  "
$ws.Range("E3").Value = "  -0.86%  "
$ws.Range("E4").Value = "  +0.32%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.87"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  -2.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5981"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  -4.45%  "
$ws.Range("E7").Value = "  +0.22%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.06964"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = "  -5.83%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2751"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = "  -4.88%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.26"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  -6.52%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07600"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = "  -1.52%  "
$ws.Range("D12").Value = "1.836.09"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.758"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = "  -4.21%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6258"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  -6.82%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000009649"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = "  -7.05%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "78.29"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = "  -4.28%  "
$ws.Range("D17").Value = "28.824.15"
$ws.Range("E17").Value = "  -1.86%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.708"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = "  -8.99%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "221.04"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = "  -5.63%  "
$ws.Range("E20").Value = "  +0.27%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.53"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  -6.24%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.864"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  -5.94%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.003"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  +0.32%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "155.51"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  -0.91%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "7.962"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  -6.10%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1290"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = "  -4.14%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.52"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  -4.62%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06501"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  -10.54%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.451"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  -3.02%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.435"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = "  -3.19%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.837"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = "  -4.85%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.759"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  -6.91%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.093"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  -5.98%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.720"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = "  -5.29%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6442"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = "  -9.27%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.539"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = "  -1.67%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.733"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = "  -2.15%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01743"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  -5.23%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.524"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  -3.85%  "
$ws.Range("D40").Value = "1.171.90"
$ws.Range("E40").Value = "  -5.00%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8929"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  -6.63%  "
$ws.Range("D43").Value = "1.981.92"
$ws.Range("E43").Value = "  -0.55%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.42"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "62.04"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  -5.02%  "
$ws.Range("E46").Value = "  -5.53%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.471"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  -5.34%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05561"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = "  -1.83%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.586"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  -6.74%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4550"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  -0.58%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3641"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  -6.26%  "
